$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full refreshed roster table (header row A1:C1 is untouched).
# Row layout: Player, Position(s), Team
$rows = @(
    @("Scottie Barnes",  "PG,SG,SF,PF", "Toronto Raptors"),
    @("Tyler Herro",     "PG,SG",       "Miami Heat"),
    @("Josh Giddey",     "PG,SG,SF",    "Chicago Bulls"),
    @("Mikal Bridges",   "SG,SF,PF",    "New York Knicks"),
    @("Miles Bridges",   "SF,PF",       "Charlotte Hornets"),
    @("Kyle Kuzma",      "PF",          "Washington Wizards"),
    @("Naz Reid",        "PF,C",        "Minnesota Timberwolves"),
    @("Nikola Vucevic",  "PF,C",        "Chicago Bulls"),
    @("Goga Bitadze",    "C",           "Orlando Magic"),
    @("Brook Lopez",     "C",           "Milwaukee Bucks"),
    @("Nikola Jovic",    "PF,C",        "Miami Heat"),
    @("DeMar DeRozan",   "SF,PF",       "Sacramento Kings"),
    @("Jeremy Sochan",   "SF,PF",       "San Antonio Spurs"),
    @("Ausar Thompson",  "SF,PF",       "Detroit Pistons"),
    @("Evan Mobley",     "PF,C",        "Cleveland Cavaliers"),
    @("Luka Doncic",     "PG,SG",       "Dallas Mavericks"),
    @("Ja Morant",       "PG",          "Memphis Grizzlies"),
    @("De'Aaron Fox",    "PG",          "Sacramento Kings")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
